# Update the "想去人数" (interest count, column F) values on the "展览"
# and "全部类型" worksheets to reflect the latest scrape.

$wb = $excel.ActiveWorkbook

# Map of worksheet name -> row number -> new value for column F.
$updates = @{
    "展览" = @{
        2  = 10085
        13 = 3137
        21 = 551
        22 = 52
        23 = 235
        26 = 231
        30 = 358
        35 = 24
        36 = 314
        37 = 1655
        39 = 416
        42 = 938
        44 = 347
    }
    "全部类型" = @{
        2  = 10085
        14 = 3137
        21 = 551
        22 = 52
        23 = 235
        26 = 231
        30 = 358
        38 = 24
        40 = 314
        41 = 1655
        44 = 416
        47 = 938
        49 = 347
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Range("F$row").Value = $rows[$row]
    }
}
